# Applies the Wed May 17 16:19:01 UTC 2023 cryptos-list refresh.
# Updates per-row Price (D) and Volume(1h) (E) figures, and for the two
# rows whose ranking swapped places, also updates Coin (B) and Link (C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores prices as literal strings straight from the feed (note
# the European-style thousands separator in values such as "26.763.74").
# The refreshed prices below that would otherwise parse as plain numbers
# (e.g. "1.001") need their cells pre-formatted as Text so they keep
# round-tripping as literal strings instead of being coerced to numeric
# values. (Applied range-by-range; this host does not honour a single
# comma-joined multi-area Range address.)
$ws.Range("D4:D11").NumberFormat = "@"
$ws.Range("D13:D20").NumberFormat = "@"
$ws.Range("D22:D36").NumberFormat = "@"
$ws.Range("D38:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '26.763.74'
$ws.Range("E2").Value = '  -1.10%  '

# Row 3
$ws.Range("D3").Value = '1.796.00'
$ws.Range("E3").Value = '  -1.43%  '

# Row 4
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '309.10'

# Row 6
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.02%  '

# Row 7
$ws.Range("D7").Value = '0.4396'
$ws.Range("E7").Value = '  +4.15%  '

# Row 8
$ws.Range("D8").Value = '0.3669'
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("D9").Value = '0.07361'
$ws.Range("E9").Value = '  +2.17%  '

# Row 10
$ws.Range("D10").Value = '0.8541'
$ws.Range("E10").Value = '  +1.48%  '

# Row 11
$ws.Range("D11").Value = '20.61'
$ws.Range("E11").Value = '  -0.82%  '

# Row 12
$ws.Range("D12").Value = '1.815.08'
$ws.Range("E12").Value = '  -0.39%  '

# Row 13
$ws.Range("D13").Value = '6.589'
$ws.Range("E13").Value = '  -1.18%  '

# Row 14
$ws.Range("D14").Value = '91.99'
$ws.Range("E14").Value = '  +2.64%  '

# Row 15
$ws.Range("D15").Value = '0.07059'
$ws.Range("E15").Value = '  -0.18%  '

# Row 16
$ws.Range("D16").Value = '5.247'
$ws.Range("E16").Value = '  -0.70%  '

# Row 17
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  -0.08%  '

# Row 18
$ws.Range("D18").Value = '0.000008623'
$ws.Range("E18").Value = '  -1.50%  '

# Row 19
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.07%  '

# Row 20
$ws.Range("D20").Value = '14.71'
$ws.Range("E20").Value = '  -1.26%  '

# Row 21
$ws.Range("D21").Value = '26.800.42'
$ws.Range("E21").Value = '  -1.22%  '

# Row 22
$ws.Range("D22").Value = '5.131'
$ws.Range("E22").Value = '  +0.10%  '

# Row 23
$ws.Range("D23").Value = '10.77'
$ws.Range("E23").Value = '  -0.72%  '

# Row 24
$ws.Range("D24").Value = '1.972'
$ws.Range("E24").Value = '  -1.11%  '

# Row 25
$ws.Range("D25").Value = '151.61'
$ws.Range("E25").Value = '  -0.22%  '

# Row 26
$ws.Range("D26").Value = '2.187'

# Row 27
$ws.Range("D27").Value = '18.35'
$ws.Range("E27").Value = '  +0.69%  '

# Row 28
$ws.Range("D28").Value = '5.169'
$ws.Range("E28").Value = '  -2.21%  '

# Row 29
$ws.Range("D29").Value = '117.09'
$ws.Range("E29").Value = '  +0.46%  '

# Row 30
$ws.Range("D30").Value = '0.08778'
$ws.Range("E30").Value = '  +0.49%  '

# Row 31
$ws.Range("B31").Value = 'ARBITRUM'
$ws.Range("C31").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D31").Value = '1.153'
$ws.Range("E31").Value = '  -2.48%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '0.7334'
$ws.Range("E32").Value = '  -0.81%  '

# Row 33
$ws.Range("D33").Value = '4.428'
$ws.Range("E33").Value = '  +0.14%  '

# Row 34
$ws.Range("D34").Value = '2.882'
$ws.Range("E34").Value = '  -2.16%  '

# Row 35
$ws.Range("D35").Value = '0.9992'
$ws.Range("E35").Value = '  -0.08%  '

# Row 36
$ws.Range("D36").Value = '1.085'
$ws.Range("E36").Value = '  -0.41%  '

# Row 37
$ws.Range("E37").Value = '  +0.11%  '

# Row 38
$ws.Range("D38").Value = '0.05156'
$ws.Range("E38").Value = '  -1.72%  '

# Row 39
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '7.041'
$ws.Range("E39").Value = '  -3.99%  '

# Row 40
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.5188'
$ws.Range("E40").Value = '  +2.75%  '

# Row 41
$ws.Range("D41").Value = '2.803'
$ws.Range("E41").Value = '  -2.26%  '

# Row 42
$ws.Range("D42").Value = '0.1670'
$ws.Range("E42").Value = '  -1.02%  '

# Row 43
$ws.Range("D43").Value = '8.406'
$ws.Range("E43").Value = '  -2.27%  '

# Row 44
$ws.Range("D44").Value = '0.4951'
$ws.Range("E44").Value = '  +4.95%  '

# Row 45
$ws.Range("D45").Value = '1.973'
$ws.Range("E45").Value = '  +2.03%  '

# Row 46
$ws.Range("D46").Value = '10.29'
$ws.Range("E46").Value = '  -2.51%  '

# Row 47
$ws.Range("D47").Value = '104.32'
$ws.Range("E47").Value = '  -1.98%  '

# Row 48
$ws.Range("D48").Value = '0.9991'
$ws.Range("E48").Value = '  -0.06%  '

# Row 49
$ws.Range("D49").Value = '1.654'
$ws.Range("E49").Value = '  +0.14%  '

# Row 50
$ws.Range("D50").Value = '0.06305'
$ws.Range("E50").Value = '  -0.57%  '

# Row 51
$ws.Range("D51").Value = '0.9132'
$ws.Range("E51").Value = '  +1.31%  '

